# Updated cryptos list - apply new Price (D) and Volume(1h) (E) values per row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.665.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.123.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.39%  "
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5285"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4538"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.02"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09094"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.183"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.126.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.855"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.089"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001178"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.014"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06725"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.012"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.341"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.728.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.395"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.376.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.567"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.200"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1081"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.657"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.372"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.026"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.079"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02661"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06892"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2322"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6952"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.281"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6479"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.336"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.770"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000365"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.82%  "
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07308"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.34%  "
